# CGR3 greenhouse growth traits.xlsx
#
# Commit: "Standardize some naming and formatting"
#
# 1. On "greenhouse growth" (sheet 1), column A held the text genotype
#    labels "36625-8" / "36625-10" / "36625-14" for rows 12-21 / 22-31 /
#    32-41. Replace them with the plain numeric genotype codes 8 / 10 / 14
#    (the genotype labels are standardized to bare numbers, matching the
#    "WT" rows which were already left as text).
# 2. Same standardization on "SPAD and LMA" (sheet 2), rows 2-7 / 8-13 /
#    14-19.
# 3. The author's last save left the selection on A41 on the first sheet
#    (no longer the active tab) and moved to the second sheet with A18
#    selected, which is now the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "greenhouse growth"
$ws2 = $wb.Worksheets.Item(2)   # "SPAD and LMA"

# --- greenhouse growth: genotype column A, rows 12-41 -----------------
for ($r = 12; $r -le 21; $r++) { $ws1.Cells.Item($r, 1).Value = 8 }
for ($r = 22; $r -le 31; $r++) { $ws1.Cells.Item($r, 1).Value = 10 }
for ($r = 32; $r -le 41; $r++) { $ws1.Cells.Item($r, 1).Value = 14 }

# --- SPAD and LMA: genotype column A, rows 2-19 ------------------------
for ($r = 2;  $r -le 7;  $r++) { $ws2.Cells.Item($r, 1).Value = 8 }
for ($r = 8;  $r -le 13; $r++) { $ws2.Cells.Item($r, 1).Value = 10 }
for ($r = 14; $r -le 19; $r++) { $ws2.Cells.Item($r, 1).Value = 14 }

# --- view state: selection moved to A41 on sheet 1, then sheet 2 is
#     activated with A18 selected -------------------------------------
[void]$ws1.Range("A41").Select()
[void]$ws2.Range("A18").Select()
[void]$ws2.Activate()
